$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column L (shifts old L..AI to M..AJ), matching "value (calculée)" addition
$ws.Columns("L:L").Insert()

# Row 1
$ws.Range("A1").Value2 = "line_number"
$ws.Range("B1").Value2 = "source_type"
$ws.Range("C1").Value2 = "year"
$ws.Range("D1").Value2 = "customs_region"
$ws.Range("E1").Value2 = "customs_office"
$ws.Range("F1").Value2 = "partner"
$ws.Range("G1").Value2 = "export_import"
$ws.Range("H1").Value2 = "product"
$ws.Range("I1").Value2 = "origin"
$ws.Range("J1").Value2 = "width_in_line"
$ws.Range("K1").Value2 = "value (annoncée)"
$ws.Range("L1").Value2 = "value (calculée)"
$ws.Range("M1").Value2 = "value_part_of_bundle"
$ws.Range("N1").Value2 = "quantity"
$ws.Range("O1").Value2 = "quantity_unit"
$ws.Range("P1").Value2 = "value_per_unit"
$ws.Range("Q1").Value2 = "filepath"
$ws.Range("R1").Value2 = "source"
$ws.Range("S1").Value2 = "sheet"
$ws.Range("T1").Value2 = "value_total"
$ws.Range("U1").Value2 = "value_sub_total_1"
$ws.Range("V1").Value2 = "value_sub_total_2"
$ws.Range("W1").Value2 = "value_sub_total_3"
$ws.Range("X1").Value2 = "data_collector"
$ws.Range("Y1").Value2 = "unverified"
$ws.Range("Z1").Value2 = "remarks"
$ws.Range("AA1").Value2 = "value_minus_unit_val_x_qty"
$ws.Range("AB1").Value2 = "absurd_observation"
$ws.Range("AC1").Value2 = "trade_deficit"
$ws.Range("AD1").Value2 = "trade_surplus"
$ws.Range("AE1").Value2 = "duty_quantity"
$ws.Range("AF1").Value2 = "duty_quantity_unit"
$ws.Range("AG1").Value2 = "duty_by_unit"
$ws.Range("AH1").Value2 = "duty_total"
$ws.Range("AI1").Value2 = "duty_part_of_bundle"
$ws.Range("AJ1").Value2 = "duty_remarks"

# Row 2
$ws.Range("A2").Value2 = 1
$ws.Range("B2").Value2 = "Local"
$ws.Range("C2").Value2 = 1749
$ws.Range("D2").Value2 = "Marseille"
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value2 = "Imports"
$ws.Range("H2").Value2 = "Balais de palme"
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value2 = 22
$ws.Range("M2").Value2 = 0
$ws.Range("N2").Value2 = 25
$ws.Range("O2").Value2 = "douzaine"
$ws.Range("Q2").Value2 = "Local/Marseille/Archives_de_la_CCI_de_Marseille-I32/Marseille – Imports – 1749.csv`t"
$ws.Range("R2").Value2 = "Archives de la CCI de Marseille - I 32"
$ws.Range("S2").Value2 = 1
$ws.Range("T2").ClearContents()
$ws.Range("U2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("W2").ClearContents()
$ws.Range("X2").Value2 = "Guillaume Daudin"
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("AC2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()
$ws.Range("AF2").ClearContents()
$ws.Range("AG2").ClearContents()
$ws.Range("AH2").ClearContents()
$ws.Range("AI2").ClearContents()
$ws.Range("AJ2").ClearContents()

# Row 3
$ws.Range("A3").Value2 = 2
$ws.Range("B3").Value2 = "Local"
$ws.Range("C3").Value2 = 1749
$ws.Range("D3").Value2 = "Marseille"
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value2 = "Imports"
$ws.Range("H3").Value2 = "Bœuf salé"
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("M3").Value2 = 0
$ws.Range("N3").Value2 = 33700
$ws.Range("O3").Value2 = "livres"
$ws.Range("Q3").Value2 = "Local/Marseille/Archives_de_la_CCI_de_Marseille-I32/Marseille – Imports – 1749.csv`t"
$ws.Range("R3").Value2 = "Archives de la CCI de Marseille - I 32"
$ws.Range("S3").Value2 = 1
$ws.Range("T3").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("V3").ClearContents()
$ws.Range("W3").ClearContents()
$ws.Range("X3").Value2 = "Guillaume Daudin"
$ws.Range("Y3").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
$ws.Range("AC3").ClearContents()
$ws.Range("AD3").ClearContents()
$ws.Range("AE3").ClearContents()
$ws.Range("AF3").ClearContents()
$ws.Range("AG3").ClearContents()
$ws.Range("AH3").ClearContents()
$ws.Range("AI3").ClearContents()
$ws.Range("AJ3").ClearContents()

# Row 4
$ws.Range("A4").Value2 = 3
$ws.Range("B4").Value2 = "Local"
$ws.Range("C4").Value2 = 1749
$ws.Range("D4").Value2 = "Marseille"
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value2 = "Imports"
$ws.Range("H4").Value2 = "Bierre"
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value2 = 75
$ws.Range("M4").Value2 = 0
$ws.Range("N4").Value2 = 1500
$ws.Range("O4").Value2 = "livres"
$ws.Range("Q4").Value2 = "Local/Marseille/Archives_de_la_CCI_de_Marseille-I32/Marseille – Imports – 1749.csv`t"
$ws.Range("R4").Value2 = "Archives de la CCI de Marseille - I 32"
$ws.Range("S4").Value2 = 1
$ws.Range("T4").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("V4").ClearContents()
$ws.Range("W4").ClearContents()
$ws.Range("X4").Value2 = "Guillaume Daudin"
$ws.Range("Y4").ClearContents()
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
$ws.Range("AC4").ClearContents()
$ws.Range("AD4").ClearContents()
$ws.Range("AE4").ClearContents()
$ws.Range("AF4").ClearContents()
$ws.Range("AG4").ClearContents()
$ws.Range("AH4").ClearContents()
$ws.Range("AI4").ClearContents()
$ws.Range("AJ4").ClearContents()

# Row 5
$ws.Range("A5").Value2 = 4
$ws.Range("B5").Value2 = "Local"
$ws.Range("C5").Value2 = 1749
$ws.Range("D5").Value2 = "Marseille"
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value2 = "Imports"
$ws.Range("H5").Value2 = "Bled froment"
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value2 = 9242650
$ws.Range("M5").Value2 = 0
$ws.Range("N5").Value2 = 369706
$ws.Range("O5").Value2 = "charges"
$ws.Range("P5").Value2 = 25
$ws.Range("Q5").Value2 = "Local/Marseille/Archives_de_la_CCI_de_Marseille-I32/Marseille – Imports – 1749.csv`t"
$ws.Range("R5").Value2 = "Archives de la CCI de Marseille - I 32"
$ws.Range("S5").Value2 = 1
$ws.Range("T5").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("V5").ClearContents()
$ws.Range("W5").ClearContents()
$ws.Range("X5").Value2 = "Guillaume Daudin"
$ws.Range("Y5").ClearContents()
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()
$ws.Range("AC5").ClearContents()
$ws.Range("AD5").ClearContents()
$ws.Range("AE5").ClearContents()
$ws.Range("AF5").ClearContents()
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AI5").ClearContents()
$ws.Range("AJ5").ClearContents()

# Row 6
$ws.Range("A6").Value2 = 5
$ws.Range("B6").Value2 = "Local"
$ws.Range("C6").Value2 = 1749
$ws.Range("D6").Value2 = "Marseille"
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("G6").Value2 = "Imports"
$ws.Range("H6").Value2 = "Bois à brûler"
$ws.Range("I6").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("K6").Value2 = 270
$ws.Range("M6").Value2 = 0
$ws.Range("N6").Value2 = 30000
$ws.Range("O6").Value2 = "livres"
$ws.Range("Q6").Value2 = "Local/Marseille/Archives_de_la_CCI_de_Marseille-I32/Marseille – Imports – 1749.csv`t"
$ws.Range("R6").Value2 = "Archives de la CCI de Marseille - I 32"
$ws.Range("S6").Value2 = 1
$ws.Range("T6").ClearContents()
$ws.Range("U6").ClearContents()
$ws.Range("V6").ClearContents()
$ws.Range("W6").ClearContents()
$ws.Range("X6").Value2 = "Guillaume Daudin"
$ws.Range("Y6").ClearContents()
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()
$ws.Range("AC6").ClearContents()
$ws.Range("AD6").ClearContents()
$ws.Range("AE6").ClearContents()
$ws.Range("AF6").ClearContents()
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").ClearContents()
$ws.Range("AJ6").ClearContents()

# --- Formulas ---
# column L: value (calculée) = quantity * value_per_unit, shared across L3:L6
$ws.Range("L2").Formula = "=N2*P2"
$ws.Range("L3:L6").Formula = "=N3*P3"

# column P: value_per_unit expressed as fractions
$ws.Range("P2").Formula = "=18/20"
$ws.Range("P3").Formula = "=7/20"
$ws.Range("P4").Formula = "=1/20"
$ws.Range("P6").Formula = "=18/20/100"

# K3: literal formula =11795
$ws.Range("K3").Formula = "=11795"

# column AA: value_minus_unit_val_x_qty = value (annoncée) - value (calculée), shared across AA3:AA6
$ws.Range("AA2").Formula = "=K2-L2"
$ws.Range("AA3:AA6").Formula = "=K3-L3"

# --- Styling: column R (source) uses Verdana 10pt font ---
$ws.Range("R2:R6").Font.Name = "Verdana"
$ws.Range("R2:R6").Font.Size = 10

# --- View state: selection and visible area ---
$ws.Range("P5").Select()
